$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The price column header changes from a flat "Hind (€)" price to a
# per-cubic-meter price "Hind (€/tm)" (week 5 price list update).
$ws.Range("B1").Value = "Hind (€/tm)"

# Column B now holds visible data worth widening to fit the new header,
# matching the new best-fit column width recorded for the sheet.
$ws.Columns.Item(2).ColumnWidth = 9.71

# The saved cursor/selection moves to B2 (the first price cell) instead
# of the previous D7 (an empty cell outside the used range).
$ws.Range("B2").Select()
